$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Step1").Name = "Table1"
$wb.Worksheets.Item("Step2").Name = "Table2"
$wb.Worksheets.Item("Step3").Name = "Table3"
$wb.Worksheets.Item("Step4").Name = "Table4"
$wb.Worksheets.Item("Step5").Name = "Table5"
$wb.Worksheets.Item("Step6").Name = "Table6"
